# CS361 Chronotimer Test Plan Template
# Author's commit: "Updating spreadsheet with my part of the lab" - fills in
# the rest of the "Test Cases Test Plan ID 5" sheet (TC05.02's remaining
# columns plus TC05.03-TC05.07) and the matching "Test Plan" summary block
# (Resources/Testers/Scheduling/... for test plan #5).
#
# Cell values are written in the same order the author apparently entered
# them (reconstructed from the shared-string table order in the target
# file) so the rebuilt shared-string table lines up exactly.

$wb = $excel.ActiveWorkbook

$tc5 = $wb.Worksheets.Item("Test Cases Test Plan ID 5")
$planSheet = $wb.Worksheets.Item("Test Plan")

# ---------------------------------------------------------------------------
# Sheet "Test Cases Test Plan ID 5": Test Case IDs + descriptions, scattered
# across rows 2-8 as they were originally authored.
# ---------------------------------------------------------------------------
$tc5.Range("A4").Value = "TC05.03"
$tc5.Range("A5").Value = "TC05.04"
$tc5.Range("A6").Value = "TC05.05"
$tc5.Range("B4").Value = "newrun while system on and current run exists"
$tc5.Range("B6").Value = "newrun with PARIND event type selected and no current run"
$tc5.Range("B7").Value = "newrun with GRP event type selected and no current run"
$tc5.Range("A7").Value = "TC05.06"
$tc5.Range("A8").Value = "TC05.07"
$tc5.Range("B3").Value = "newrun while system on and no current run and no run type specified"
$tc5.Range("B5").Value = "newrun with IND event type selected and no current run"
$tc5.Range("B8").Value = "newrun with PARGRP event type selected and no current run"
$tc5.Range("B2").Value = "newrun While System OFF"

# Features being tested / Input values - columns C and D. C2/D2/C5 carry the
# leftover "blank but bordered" style from the template, so drop that style
# before writing real text into them (matches the target file, which shows
# these cells with no explicit style once filled in).
$tc5.Range("C2").ClearFormats()
$tc5.Range("C2").Value = "while power off - <time> newrun"
$tc5.Range("D2").ClearFormats()
$tc5.Range("D2").Value = "currentRun == null"
$tc5.Range("D4").Value = "IllegalStateException"
$tc5.Range("C4").Value = "while power on and currentRun!=null - <time> newrun"
$tc5.Range("D3").Value = "currentRun == new Run"
$tc5.Range("C3").Value = "while power on and currentRun==null - <time> newrun"
$tc5.Range("C5").ClearFormats()
$tc5.Range("C5").Value = 'while power on, currentRun==null, and eventType="IND" - <time> newrun'
$tc5.Range("D5").Value = "currentRun == new IND()"
$tc5.Range("C6").Value = 'while power on, currentRun==null, and eventType="PARIND" - <time> newrun'
$tc5.Range("D6").Value = "currentRun == new PARIND()"
$tc5.Range("C7").Value = 'while power on, currentRun==null, and eventType="GRP" - <time> newrun'
$tc5.Range("D7").Value = "currentRun == new GRP()"
$tc5.Range("C8").Value = 'while power on, currentRun==null, and eventType="PARGRP" - <time> newrun'
$tc5.Range("D8").Value = "currentRun == new PARGRP()"

# ---------------------------------------------------------------------------
# Sheet "Test Plan": test plan #5's block (rows 55-62), column B.
# ---------------------------------------------------------------------------
$planSheet.Range("B55").Value = "Time on 1 Computer"
$planSheet.Range("B57").Value = "TBD"
$planSheet.Range("B58").Value = "N/A"
$planSheet.Range("B59").Value = "System off"
$planSheet.Range("B60").Value = "newrun"
$planSheet.Range("B61").Value = "power"
$planSheet.Range("B62").Value = "event"

# ---------------------------------------------------------------------------
# Back to the test-case sheet: Expected Output/Actual/Author/Tester columns
# are still "TBD" placeholders for every row.
# ---------------------------------------------------------------------------
foreach ($row in 2..8) {
    foreach ($col in "E", "F", "G", "H") {
        $tc5.Range("$col$row").Value = "TBD"
    }
}

# Widen the new columns so the long descriptions are readable, then restore
# the selection/active sheet the workbook was left on.
$tc5.Columns("A:D").AutoFit()

# Last cell the author touched on the "Testers:" row.
$planSheet.Range("B56").Value = "1 Member of Team TBD"

$tc5.Range("D12").Select()
$planSheet.Activate()
$planSheet.Range("B58").Select()
